# Update the two-digit multiplication problems in the table with a new
# set of operands, matching the data-refresh described in the commit
# message ("Update master to output generated at c986bee").
#
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#              Format, ReplaceWith, Replace)
#   MatchCase      = $true  -> exact text match
#   Forward        = $true
#   Wrap           = 1      -> wdFindContinue (search whole story once)
#   ReplaceWith    = new problem text
#   Replace        = 2      -> wdReplaceAll (each old string appears once)
$d = $word.ActiveDocument

$d.Content.Find.Execute("12×19=", $true, $false, $false, $false, $false, $true, 1, $false, "95×27=", 2) | Out-Null
$d.Content.Find.Execute("39×21=", $true, $false, $false, $false, $false, $true, 1, $false, "63×87=", 2) | Out-Null
$d.Content.Find.Execute("60×60=", $true, $false, $false, $false, $false, $true, 1, $false, "56×28=", 2) | Out-Null
$d.Content.Find.Execute("68×12=", $true, $false, $false, $false, $false, $true, 1, $false, "78×90=", 2) | Out-Null
$d.Content.Find.Execute("85×78=", $true, $false, $false, $false, $false, $true, 1, $false, "75×51=", 2) | Out-Null
$d.Content.Find.Execute("26×23=", $true, $false, $false, $false, $false, $true, 1, $false, "26×66=", 2) | Out-Null
$d.Content.Find.Execute("86×21=", $true, $false, $false, $false, $false, $true, 1, $false, "95×23=", 2) | Out-Null
$d.Content.Find.Execute("33×66=", $true, $false, $false, $false, $false, $true, 1, $false, "20×80=", 2) | Out-Null
$d.Content.Find.Execute("23×83=", $true, $false, $false, $false, $false, $true, 1, $false, "23×76=", 2) | Out-Null
$d.Content.Find.Execute("51×14=", $true, $false, $false, $false, $false, $true, 1, $false, "54×81=", 2) | Out-Null
$d.Content.Find.Execute("98×45=", $true, $false, $false, $false, $false, $true, 1, $false, "15×51=", 2) | Out-Null
$d.Content.Find.Execute("61×55=", $true, $false, $false, $false, $false, $true, 1, $false, "31×24=", 2) | Out-Null
$d.Content.Find.Execute("21×71=", $true, $false, $false, $false, $false, $true, 1, $false, "13×77=", 2) | Out-Null
$d.Content.Find.Execute("59×23=", $true, $false, $false, $false, $false, $true, 1, $false, "30×79=", 2) | Out-Null
$d.Content.Find.Execute("75×47=", $true, $false, $false, $false, $false, $true, 1, $false, "69×71=", 2) | Out-Null
$d.Content.Find.Execute("60×18=", $true, $false, $false, $false, $false, $true, 1, $false, "36×78=", 2) | Out-Null
$d.Content.Find.Execute("38×97=", $true, $false, $false, $false, $false, $true, 1, $false, "92×92=", 2) | Out-Null
$d.Content.Find.Execute("26×42=", $true, $false, $false, $false, $false, $true, 1, $false, "94×22=", 2) | Out-Null
$d.Content.Find.Execute("88×25=", $true, $false, $false, $false, $false, $true, 1, $false, "17×14=", 2) | Out-Null
$d.Content.Find.Execute("75×52=", $true, $false, $false, $false, $false, $true, 1, $false, "20×31=", 2) | Out-Null
$d.Content.Find.Execute("48×81=", $true, $false, $false, $false, $false, $true, 1, $false, "63×50=", 2) | Out-Null
$d.Content.Find.Execute("12×82=", $true, $false, $false, $false, $false, $true, 1, $false, "81×84=", 2) | Out-Null
$d.Content.Find.Execute("91×33=", $true, $false, $false, $false, $false, $true, 1, $false, "49×69=", 2) | Out-Null
$d.Content.Find.Execute("88×72=", $true, $false, $false, $false, $false, $true, 1, $false, "39×13=", 2) | Out-Null
$d.Content.Find.Execute("65×93=", $true, $false, $false, $false, $false, $true, 1, $false, "59×96=", 2) | Out-Null
